# "start adding in asbuilt data"
#
# 1. Cell K9 held a stray "X" marker - blank it out to a single space.
# 2. Turn on an AutoFilter over the existing table (A1:J67) that filters
#    column D ("tier") down to just the "detector-simulated" rows. That
#    hides every data row except row 2 (the one row whose tier already
#    reads "detector-simulated").
# 3. The saved selection moves to K1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the leftover "X" in K9 to a blank space.
$ws.Range("K9").Value = " "

# 2. Apply the AutoFilter: Field 4 = column D ("tier"), filtered to the
#    single value "detector-simulated" (xlFilterValues = 7).
$ws.Range("A1:J67").AutoFilter(4, @("detector-simulated"), 7)

# Row 2 already matches "detector-simulated" and should stay visible;
# everything else in the table (rows 3-67) ends up hidden by the filter.
$ws.Rows.Item(2).Hidden = $false

# 3. Leave the selection on K1, matching the saved sheet view.
$ws.Range("K1").Select()
